$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore cell C10 ("Integer min" for rule R30) from 18 back to 1,
# as per the source revision being restored.
$ws.Range("C10").Value = 1
